$d = $word.ActiveDocument

$d.Content.Find.Execute("959÷7=137, 0", $true, $false, $false, $false, $false, $true, 1, $false, "563÷7=80, 3", 2) | Out-Null
$d.Content.Find.Execute("764÷5=152, 4", $true, $false, $false, $false, $false, $true, 1, $false, "378÷3=126, 0", 2) | Out-Null
$d.Content.Find.Execute("408÷4=102, 0", $true, $false, $false, $false, $false, $true, 1, $false, "297÷3=99, 0", 2) | Out-Null
$d.Content.Find.Execute("363÷2=181, 1", $true, $false, $false, $false, $false, $true, 1, $false, "402÷6=67, 0", 2) | Out-Null
$d.Content.Find.Execute("816÷2=408, 0", $true, $false, $false, $false, $false, $true, 1, $false, "480÷6=80, 0", 2) | Out-Null
$d.Content.Find.Execute("573÷6=95, 3", $true, $false, $false, $false, $false, $true, 1, $false, "526÷2=263, 0", 2) | Out-Null
$d.Content.Find.Execute("957÷9=106, 3", $true, $false, $false, $false, $false, $true, 1, $false, "739÷4=184, 3", 2) | Out-Null
$d.Content.Find.Execute("697÷2=348, 1", $true, $false, $false, $false, $false, $true, 1, $false, "848÷8=106, 0", 2) | Out-Null
$d.Content.Find.Execute("578÷2=289, 0", $true, $false, $false, $false, $false, $true, 1, $false, "425÷5=85, 0", 2) | Out-Null
$d.Content.Find.Execute("984÷8=123, 0", $true, $false, $false, $false, $false, $true, 1, $false, "110÷6=18, 2", 2) | Out-Null
$d.Content.Find.Execute("782÷8=97, 6", $true, $false, $false, $false, $false, $true, 1, $false, "321÷6=53, 3", 2) | Out-Null
$d.Content.Find.Execute("893÷5=178, 3", $true, $false, $false, $false, $false, $true, 1, $false, "830÷8=103, 6", 2) | Out-Null
$d.Content.Find.Execute("260÷7=37, 1", $true, $false, $false, $false, $false, $true, 1, $false, "591÷9=65, 6", 2) | Out-Null
$d.Content.Find.Execute("330÷6=55, 0", $true, $false, $false, $false, $false, $true, 1, $false, "110÷4=27, 2", 2) | Out-Null
$d.Content.Find.Execute("565÷4=141, 1", $true, $false, $false, $false, $false, $true, 1, $false, "684÷3=228, 0", 2) | Out-Null
$d.Content.Find.Execute("796÷9=88, 4", $true, $false, $false, $false, $false, $true, 1, $false, "236÷6=39, 2", 2) | Out-Null
$d.Content.Find.Execute("857÷7=122, 3", $true, $false, $false, $false, $false, $true, 1, $false, "513÷8=64, 1", 2) | Out-Null
$d.Content.Find.Execute("265÷3=88, 1", $true, $false, $false, $false, $false, $true, 1, $false, "911÷5=182, 1", 2) | Out-Null
$d.Content.Find.Execute("759÷3=253, 0", $true, $false, $false, $false, $false, $true, 1, $false, "885÷2=442, 1", 2) | Out-Null
$d.Content.Find.Execute("694÷9=77, 1", $true, $false, $false, $false, $false, $true, 1, $false, "424÷8=53, 0", 2) | Out-Null
$d.Content.Find.Execute("377÷8=47, 1", $true, $false, $false, $false, $false, $true, 1, $false, "251÷2=125, 1", 2) | Out-Null
$d.Content.Find.Execute("797÷2=398, 1", $true, $false, $false, $false, $false, $true, 1, $false, "792÷7=113, 1", 2) | Out-Null
$d.Content.Find.Execute("464÷7=66, 2", $true, $false, $false, $false, $false, $true, 1, $false, "859÷8=107, 3", 2) | Out-Null
$d.Content.Find.Execute("347÷2=173, 1", $true, $false, $false, $false, $false, $true, 1, $false, "174÷5=34, 4", 2) | Out-Null
$d.Content.Find.Execute("940÷7=134, 2", $true, $false, $false, $false, $false, $true, 1, $false, "578÷8=72, 2", 2) | Out-Null
